$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old B1 cell carried a custom date-number-format style (xf index 1).
# That header cell becomes a plain text header ("Column2"), so drop the
# inherited formatting first -- otherwise the style sticks around on the
# cell even after the value/type changes.
$ws.Range("B1").ClearFormats()

# Shift the original single data row down to row 2, and put fresh
# "ColumnN" headers on row 1 (this mirrors the row A1:E1 -> A2:E2 move +
# new header row visible in the diff).
$ws.Range("A2").Value = "This is a string in A1"
$ws.Range("B2").Value = 44910.578125254629
$ws.Range("C2").Value = 1234
$ws.Range("D2").Value = 123456
$ws.Range("E2").Value = 1.234

$ws.Range("A1").Value = "Column1"
$ws.Range("B1").Value = "Column2"
$ws.Range("C1").Value = "Column3"
$ws.Range("D1").Value = "Column4"
$ws.Range("E1").Value = "Column5"

# Columns C:E lose their old bestFit widths once they're part of the new
# table and take on the table's default column width.
$ws.Range("C1:E1").ColumnWidth = 9.666666666666666

# Turn A1:E2 into an actual Excel table ("Table1") with a light style.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:E2"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight1"

# Match the final selection left behind in the saved file.
[void]$ws.Range("F8").Select()
